$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column prices are stored as plain text in the workbook (e.g. '69.445.16',
# '1.81'), not numbers, so values that look numeric are written with a leading
# apostrophe to keep them as text, then the cell style is reset to 'Normal' so
# no stray number-format style gets attached.

$ws.Range("D2").Value = "69.456.05"
$ws.Range("E2").Value = "  -2.65%  "

$ws.Range("D3").Value = "3.692.09"
$ws.Range("E3").Value = "  -3.13%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'688.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "

$ws.Range("D6").Value = "'162.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.30%  "

$ws.Range("D7").Value = "3.692.43"
$ws.Range("E7").Value = "  -3.11%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -5.47%  "

$ws.Range("E10").Value = "  -8.07%  "

$ws.Range("D11").Value = "'7.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.71%  "

$ws.Range("E12").Value = "  -8.36%  "

$ws.Range("E13").Value = "  -5.73%  "

$ws.Range("D14").Value = "4.316.29"
$ws.Range("E14").Value = "  -3.08%  "

$ws.Range("D15").Value = "'33.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.08%  "

$ws.Range("D16").Value = "3.688.06"
$ws.Range("E16").Value = "  -2.88%  "

$ws.Range("D17").Value = "69.502.31"
$ws.Range("E17").Value = "  -2.62%  "

$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("E19").Value = "  -8.23%  "

$ws.Range("D20").Value = "'6.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.42%  "

$ws.Range("D21").Value = "'476.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.33%  "

$ws.Range("D22").Value = "'9.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.14%  "

$ws.Range("D23").Value = "'0.658"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.87%  "

$ws.Range("D24").Value = "'79.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.90%  "

$ws.Range("D25").Value = "3.838.21"
$ws.Range("E25").Value = "  -3.03%  "

$ws.Range("E26").Value = "  -9.27%  "

$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("E28").Value = "  -9.94%  "

$ws.Range("D29").Value = "'9.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.41%  "

$ws.Range("D30").Value = "'1.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.20%  "

$ws.Range("E31").Value = "  -9.91%  "

$ws.Range("D32").Value = "'6.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.77%  "

$ws.Range("D33").Value = "'2.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.03%  "

$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("D35").Value = "'26.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.73%  "

$ws.Range("E36").Value = "  -3.84%  "

$ws.Range("D37").Value = "3.657.79"

$ws.Range("D38").Value = "'8.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.94%  "

$ws.Range("D39").Value = "'6.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.58%  "

$ws.Range("D40").Value = "'2.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.66%  "

$ws.Range("D41").Value = "'0.0919"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.97%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("E44").Value = "  -6.22%  "

$ws.Range("D45").Value = "'163.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.02%  "

$ws.Range("D46").Value = "'48.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "

$ws.Range("D47").Value = "'30.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("E48").Value = "  -15.21%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "'0.000283"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.33%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.87%  "

$ws.Range("E51").Value = "  -2.52%  "
